$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNum, D(Fecha), I(Calidad), J(Volumen), K(PrecioMinimo), L(PrecioMaximo), M(PrecioPromedioPonderado), O(Origen), P(PrecioPorKg)
$data = @(
    ,@(436, 44722, 'Primera', 10600, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(437, 44722, 'Segunda', 4300, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(438, 44281, 'Primera', 4300, 110, 110, 110, 'Región Metropolitana', 110)
    ,@(439, 44623, 'Primera', 6100, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(440, 44623, 'Segunda', 3400, 80, 80, 80, 'Región Metropolitana', 80)
    ,@(441, 44313, 'Primera', 5200, 130, 130, 130, 'Región Metropolitana', 130)
    ,@(442, 44664, 'Primera', 7900, 100, 120, 110, 'Región Metropolitana', 110)
    ,@(443, 44664, 'Segunda', 3400, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(444, 44195, 'Primera', 16000, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(445, 44433, 'Primera', 5200, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(446, 44433, 'Segunda', 5200, 80, 90, 85, 'Región Metropolitana', 85)
    ,@(447, 44676, 'Primera', 5200, 120, 130, 125, 'Región Metropolitana', 125)
    ,@(448, 44292, 'Primera', 3400, 130, 130, 130, 'Región Metropolitana', 130)
    ,@(449, 44655, 'Primera', 4300, 120, 130, 125, 'Región Metropolitana', 125)
    ,@(450, 44330, 'Primera', 6100, 120, 120, 120, 'Región Metropolitana', 120)
    ,@(451, 44398, 'Primera', 5200, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(452, 44398, 'Segunda', 1600, 70, 70, 70, 'Región Metropolitana', 70)
    ,@(453, 44578, 'Primera', 4300, 100, 120, 110, 'Región Metropolitana', 110)
    ,@(454, 44578, 'Segunda', 1960, 70, 80, 75, 'Región Metropolitana', 75)
    ,@(455, 44495, 'Primera', 17500, 100, 120, 109, 'Provincia de Chacabuco', 109)
    ,@(456, 44495, 'Segunda', 6800, 80, 80, 80, 'Provincia de Chacabuco', 80)
    ,@(457, 44648, 'Primera', 4300, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(458, 44648, 'Segunda', 1600, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(459, 44221, 'Primera', 5000, 90, 100, 94, 'Provincia de Chacabuco', 94)
    ,@(460, 44580, 'Primera', 4300, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(461, 44580, 'Segunda', 2130, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(462, 44371, 'Primera', 6100, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(463, 44371, 'Segunda', 3400, 80, 80, 80, 'Región Metropolitana', 80)
    ,@(464, 44579, 'Primera', 6100, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(465, 44579, 'Segunda', 3400, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(466, 44551, 'Primera', 6100, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(467, 44551, 'Segunda', 2500, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(468, 44285, 'Primera', 4300, 110, 110, 110, 'Región Metropolitana', 110)
    ,@(469, 44314, 'Primera', 5200, 130, 130, 130, 'Región Metropolitana', 130)
    ,@(470, 44708, 'Primera', 15000, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(471, 44708, 'Segunda', 7000, 100, 100, 100, 'Región Metropolitana', 100)
    ,@(472, 44454, 'Primera', 4300, 90, 110, 100, 'Región Metropolitana', 100)
    ,@(473, 44454, 'Segunda', 1330, 70, 80, 75, 'Región Metropolitana', 75)
    ,@(474, 44160, 'Primera', 16000, 80, 100, 90, 'Provincia de Chacabuco', 90)
    ,@(475, 44554, 'Primera', 9700, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(476, 44554, 'Segunda', 3400, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(477, 44565, 'Primera', 6100, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(478, 44565, 'Segunda', 2500, 70, 70, 70, 'Región Metropolitana', 70)
    ,@(479, 44603, 'Primera', 7900, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(480, 44603, 'Segunda', 4300, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(481, 44196, 'Primera', 16000, 80, 90, 85, 'Región Metropolitana', 85)
    ,@(482, 44407, 'Primera', 6100, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(483, 44407, 'Segunda', 2500, 70, 70, 70, 'Región Metropolitana', 70)
    ,@(484, 44263, 'Primera', 4300, 110, 110, 110, 'Región Metropolitana', 110)
    ,@(485, 44187, 'Primera', 21000, 80, 100, 90, 'Provincia de Chacabuco', 90)
    ,@(486, 44609, 'Primera', 7900, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(487, 44609, 'Segunda', 2500, 70, 70, 70, 'Región Metropolitana', 70)
    ,@(488, 44529, 'Primera', 3400, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(489, 44529, 'Segunda', 1600, 80, 90, 85, 'Región Metropolitana', 85)
    ,@(490, 44321, 'Primera', 6100, 130, 130, 130, 'Región Metropolitana', 130)
    ,@(491, 44277, 'Primera', 3400, 100, 100, 100, 'Región Metropolitana', 100)
    ,@(492, 44166, 'Primera', 19000, 80, 100, 89, 'Provincia de Chacabuco', 89)
    ,@(493, 44627, 'Primera', 4300, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(494, 44354, 'Primera', 4300, 120, 130, 125, 'Región Metropolitana', 125)
    ,@(495, 44245, 'Primera', 7000, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(496, 44168, 'Primera', 22000, 80, 100, 91, 'Provincia de Chacabuco', 91)
    ,@(497, 44638, 'Primera', 6100, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(498, 44638, 'Segunda', 3400, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(499, 44522, 'Primera', 4300, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(500, 44522, 'Segunda', 2500, 80, 80, 80, 'Región Metropolitana', 80)
    ,@(501, 44699, 'Primera', 17000, 110, 120, 115, 'Provincia de Chacabuco', 115)
    ,@(502, 44699, 'Segunda', 8000, 100, 100, 100, 'Provincia de Chacabuco', 100)
    ,@(503, 44299, 'Primera', 4300, 150, 150, 150, 'Región Metropolitana', 150)
    ,@(504, 44714, 'Primera', 10600, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(505, 44714, 'Segunda', 4300, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(506, 44615, 'Primera', 5200, 100, 120, 110, 'Región Metropolitana', 110)
    ,@(507, 44615, 'Segunda', 1960, 80, 80, 80, 'Región Metropolitana', 80)
    ,@(508, 44188, 'Primera', 15000, 70, 100, 86, 'Región Metropolitana', 86)
    ,@(509, 44659, 'Primera', 7900, 100, 120, 110, 'Región Metropolitana', 110)
    ,@(510, 44659, 'Segunda', 3400, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(511, 44651, 'Primera', 9700, 110, 120, 115, 'Región Metropolitana', 115)
    ,@(512, 44651, 'Segunda', 3400, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(513, 44453, 'Primera', 5200, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(514, 44453, 'Segunda', 2500, 70, 80, 75, 'Región Metropolitana', 75)
    ,@(515, 44421, 'Primera', 5200, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(516, 44421, 'Segunda', 2500, 70, 70, 70, 'Región Metropolitana', 70)
    ,@(517, 44291, 'Primera', 4300, 130, 130, 130, 'Región Metropolitana', 130)
    ,@(518, 44323, 'Primera', 7000, 110, 110, 110, 'Región Metropolitana', 110)
    ,@(519, 44526, 'Primera', 6100, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(520, 44526, 'Segunda', 2500, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(521, 44363, 'Primera', 5200, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(522, 44363, 'Segunda', 2500, 80, 80, 80, 'Región Metropolitana', 80)
    ,@(523, 44403, 'Primera', 3400, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(524, 44403, 'Segunda', 1600, 80, 80, 80, 'Región Metropolitana', 80)
    ,@(525, 44704, 'Primera', 5000, 110, 120, 114, 'Provincia de Chacabuco', 114)
    ,@(526, 44620, 'Primera', 4300, 100, 120, 110, 'Región Metropolitana', 110)
    ,@(527, 44586, 'Primera', 7900, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(528, 44586, 'Segunda', 4300, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(529, 44601, 'Primera', 5200, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(530, 44601, 'Segunda', 2500, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(531, 44544, 'Primera', 7900, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(532, 44544, 'Segunda', 3400, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(533, 44617, 'Primera', 7900, 100, 110, 105, 'Región Metropolitana', 105)
    ,@(534, 44617, 'Segunda', 3400, 90, 90, 90, 'Región Metropolitana', 90)
    ,@(535, 44567, 'Primera', 6100, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(536, 44567, 'Segunda', 2500, 60, 70, 65, 'Región Metropolitana', 65)
    ,@(537, 44474, 'Primera', 4300, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(538, 44474, 'Segunda', 2500, 70, 80, 75, 'Región Metropolitana', 75)
    ,@(539, 44377, 'Primera', 4300, 90, 100, 95, 'Región Metropolitana', 95)
    ,@(540, 44377, 'Segunda', 1600, 80, 80, 80, 'Región Metropolitana', 80)
)

# The last two rows (539 and 540) are brand new rows appended after the
# previously-last row (538). Populate the columns that stay constant for
# every record in this sheet (A, B, C, E, F, G, H, N, Q, R) by copying them
# from row 538, before filling in the record-specific columns below.
$lastExistingRow = 538
$newRows = @(539, 540)
foreach ($newRow in $newRows) {
    $ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($lastExistingRow, 1).Value2   # A - Mercado ID
    $ws.Cells.Item($newRow, 2).Value = $ws.Cells.Item($lastExistingRow, 2).Value2   # B - Mercado
    $ws.Cells.Item($newRow, 3).Value = $ws.Cells.Item($lastExistingRow, 3).Value2   # C - Region
    $ws.Cells.Item($newRow, 5).Value = $ws.Cells.Item($lastExistingRow, 5).Value2   # E - Codreg
    $ws.Cells.Item($newRow, 6).Value = $ws.Cells.Item($lastExistingRow, 6).Value2   # F - Categoria ID
    $ws.Cells.Item($newRow, 7).Value = $ws.Cells.Item($lastExistingRow, 7).Value2   # G - Categoria
    $ws.Cells.Item($newRow, 8).Value = $ws.Cells.Item($lastExistingRow, 8).Value2   # H - Variedad
    $ws.Cells.Item($newRow, 14).Value = $ws.Cells.Item($lastExistingRow, 14).Value2 # N - Unidad de comercializacion
    $ws.Cells.Item($newRow, 17).Value = $ws.Cells.Item($lastExistingRow, 17).Value2 # Q - Kg o Unidades
    $ws.Cells.Item($newRow, 18).Value = $ws.Cells.Item($lastExistingRow, 18).Value2 # R - Clasificacion

    # Column D on existing rows uses a date/time number format (style index 2).
    # Re-apply the same style that's used on the rest of column D so the new
    # rows render consistently.
    $ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($lastExistingRow, 4).NumberFormat
}

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]    # D - Fecha
    $ws.Cells.Item($r, 9).Value = $row[2]    # I - Calidad
    $ws.Cells.Item($r, 10).Value = $row[3]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $row[4]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[5]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[6]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $row[7]   # O - Origen
    $ws.Cells.Item($r, 16).Value = $row[8]   # P - Precio $/Kg
}
